# Update the crypto price/volume snapshot on the active worksheet.
# Source values are plain text (e.g. "307.55", "0.88%") in the original
# workbook (stored as inline strings), so each cell is forced to a Text
# number format before the value is written. This keeps Excel's COM layer
# from "helpfully" re-interpreting numeric- or percent-looking text as a
# real number/percentage (which would change both the stored value and
# its formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
}

# r="2"
Set-TextValue "D2" "307.55"
Set-TextValue "E2" "0.88%"

# r="3"
Set-TextValue "D3" "36.50"
Set-TextValue "E3" "1.80%"

# r="4"
Set-TextValue "D4" "5.060"
Set-TextValue "E4" "1.64%"

# r="5"
Set-TextValue "D5" "0.08103"

# r="6"
Set-TextValue "D6" "2.031"
Set-TextValue "E6" "6.79%"

# r="7"
Set-TextValue "D7" "7.857"
Set-TextValue "E7" "-0.31%"

# r="8"
Set-TextValue "D8" "0.9287"
Set-TextValue "E8" "-0.13%"

# r="9"
Set-TextValue "D9" "0.1485"
Set-TextValue "E9" "17.16%"

# r="10"
Set-TextValue "D10" "0.1941"
Set-TextValue "E10" "1.96%"

# r="11"
Set-TextValue "D11" "0.09080"
Set-TextValue "E11" "-1.14%"

# r="12"
Set-TextValue "D12" "0.03520"
Set-TextValue "E12" "0.36%"

# r="13"
Set-TextValue "E13" "-0.36%"

# r="14"
Set-TextValue "D14" "0.001404"
Set-TextValue "E14" "-1.03%"

# r="15"
Set-TextValue "D15" "0.006125"
Set-TextValue "E15" "0.93%"

# r="16"
Set-TextValue "D16" "3.834"
Set-TextValue "E16" "6.50%"

# r="17"
Set-TextValue "D17" "4.164"
Set-TextValue "E17" "0.48%"

# r="18"
Set-TextValue "D18" "3.452"
Set-TextValue "E18" "11.18%"

# r="19"
Set-TextValue "D19" "0.3453"
Set-TextValue "E19" "0.03%"

# r="20"
Set-TextValue "D20" "0.1304"
Set-TextValue "E20" "0.68%"

# r="21"
Set-TextValue "D21" "4.833"
Set-TextValue "E21" "-7.69%"

# r="22"
Set-TextValue "D22" "0.2345"
Set-TextValue "E22" "-7.30%"

# r="23"
Set-TextValue "D23" "0.04386"
Set-TextValue "E23" "-0.52%"

# r="24"
Set-TextValue "D24" "0.001236"
Set-TextValue "E24" "0.09%"

# r="25"
Set-TextValue "D25" "0.004190"
Set-TextValue "E25" "-11.17%"

# r="27"
Set-TextValue "D27" "0.0001303"
Set-TextValue "E27" "0.22%"

# r="39"
Set-TextValue "D39" "0.02053"
Set-TextValue "E39" "4.93%"

# r="40"
Set-TextValue "D40" "0.05122"
Set-TextValue "E40" "-1.89%"

# r="41"
Set-TextValue "D41" "0.007477"
Set-TextValue "E41" "-1.10%"

# r="42"
Set-TextValue "D42" "0.01007"
Set-TextValue "E42" "-0.79%"

# r="43"
Set-TextValue "E43" "0.02%"

# r="44"
Set-TextValue "D44" "0.002125"
Set-TextValue "E44" "1.17%"

# r="45"
Set-TextValue "D45" "0.009893"
Set-TextValue "E45" "-7.12%"

# r="46"
Set-TextValue "D46" "0.00006317"
Set-TextValue "E46" "-0.56%"

# r="47"
Set-TextValue "E47" "0.13%"

# r="48"
Set-TextValue "D48" "64.85"
Set-TextValue "E48" "-0.56%"

# r="49"
Set-TextValue "D49" "0.001603"
Set-TextValue "E49" "-3.36%"

# r="50"
Set-TextValue "D50" "0.00002103"
Set-TextValue "E50" "0.13%"

# r="51"
Set-TextValue "D51" "0.0002003"
Set-TextValue "E51" "0.13%"
